# "fix waiting time and avg rate"
# Fills in the previously-zeroed "cost" (column E) figures for rows 65-105,
# clears the one row (85) that should go back to blank, tightens a batch of
# row heights that were too tall for their wrapped text, and restores the
# selection to B174.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column E value corrections (rows 65-105) ---------------------------
$costUpdates = @{
    65  = 150000
    66  = 245000
    67  = 810000
    68  = 400000
    69  = 700000
    70  = 240000
    72  = 390000
    73  = 200000
    74  = 450000
    75  = 450000
    77  = 400000
    78  = 500000
    79  = 400000
    82  = 450000
    83  = 500000
    84  = 500000
    86  = 600000
    87  = 700000
    88  = 500000
    89  = 500000
    90  = 400000
    91  = 300000
    92  = 480000
    93  = 400000
    95  = 140000
    96  = 450000
    97  = 500000
    98  = 700000
    99  = 180000
    100 = 20000
    101 = 800000
    103 = 100000
    104 = 400000
    105 = 400000
}

foreach ($row in $costUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $costUpdates[$row]
}

# Row 85's cost cell goes back to empty instead of 0
$ws.Cells.Item(85, 5).ClearContents()

# --- Row height corrections ----------------------------------------------
$rowHeights = @{
    3   = 86.4
    12  = 72
    13  = 72
    16  = 43.2
    24  = 57.6
    32  = 28.8
    33  = 57.6
    36  = 43.2
    37  = 57.6
    41  = 57.6
    56  = 28.8
    71  = 57.6
    108 = 43.2
}

foreach ($row in $rowHeights.Keys) {
    $ws.Rows.Item($row).RowHeight = $rowHeights[$row]
}

# --- Restore the saved selection -----------------------------------------
$ws.Range("B174").Select() | Out-Null
